# Update '想去人数' (interested count) figures in column F across sheets
# per 'Update gh-pages to output generated at 456a3b4'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 189
$ws.Range("F5").Value = 1038
$ws.Range("F7").Value = 2679
$ws.Range("F9").Value = 1309
$ws.Range("F10").Value = 937
$ws.Range("F11").Value = 630
$ws.Range("F12").Value = 944
$ws.Range("F13").Value = 1190
$ws.Range("F16").Value = 747
$ws.Range("F17").Value = 795
$ws.Range("F19").Value = 535
$ws.Range("F20").Value = 1139
$ws.Range("F22").Value = 650
$ws.Range("F23").Value = 613
$ws.Range("F24").Value = 232
$ws.Range("F25").Value = 319
$ws.Range("F26").Value = 316
$ws.Range("F27").Value = 698
$ws.Range("F28").Value = 592
$ws.Range("F29").Value = 5778
$ws.Range("F30").Value = 500
$ws.Range("F32").Value = 306
$ws.Range("F34").Value = 182
$ws.Range("F35").Value = 1649
$ws.Range("F37").Value = 106
$ws.Range("F38").Value = 449
$ws.Range("F40").Value = 92
$ws.Range("F42").Value = 15
$ws.Range("F43").Value = 76
$ws.Range("F45").Value = 147
$ws.Range("F46").Value = 140
$ws.Range("F47").Value = 122

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 53
$ws.Range("F12").Value = 197
$ws.Range("F13").Value = 4413

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 751

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 751
$ws.Range("F6").Value = 1038
$ws.Range("F7").Value = 2679
$ws.Range("F9").Value = 1309
$ws.Range("F10").Value = 937
$ws.Range("F11").Value = 630
$ws.Range("F12").Value = 944
$ws.Range("F13").Value = 1190
$ws.Range("F17").Value = 747
$ws.Range("F19").Value = 795
$ws.Range("F21").Value = 535
$ws.Range("F22").Value = 1139
$ws.Range("F24").Value = 53
$ws.Range("F25").Value = 650
$ws.Range("F26").Value = 613
$ws.Range("F27").Value = 232
$ws.Range("F28").Value = 319
$ws.Range("F29").Value = 316
$ws.Range("F30").Value = 592
$ws.Range("F31").Value = 5778
$ws.Range("F32").Value = 197
$ws.Range("F33").Value = 501
$ws.Range("F36").Value = 182
$ws.Range("F37").Value = 1649
$ws.Range("F39").Value = 449
$ws.Range("F42").Value = 92
$ws.Range("F44").Value = 15
$ws.Range("F45").Value = 76
$ws.Range("F46").Value = 147
$ws.Range("F48").Value = 122

